# Add a new sheet "04-20-2022" (a copy of the "Daily Attendance Template")
# positioned between "Key" and "Daily Attendance Template", pre-populated
# with some existing attendance data, so that staffRowNum initializes from
# the rows that already exist in the sheet instead of always starting at 1.

$wb = $excel.ActiveWorkbook
$templateSheet = $wb.Worksheets.Item("Daily Attendance Template")

# Copy the template sheet, placing the copy immediately before the
# original template (i.e. right after "Key").
$templateSheet.Copy($templateSheet, $null)
$ws = $wb.ActiveSheet
$ws.Name = "04-20-2022"
$ws.Unprotect()

# Existing attendance rows already logged for this day.
$ws.Range("A2").Value = "Bunk 1"
$ws.Range("B2").Value = "Staff Member 1"
$ws.Range("C2").Value = "Staff Member 1 ID"
$ws.Range("D2").Value = "10:12 PM"
$ws.Range("E2").Value = "Leaving Camp"

$ws.Range("A3").Value = "Bunk 1"
$ws.Range("B3").Value = "Staff Member 1"
$ws.Range("C3").Value = "Staff Member 1 ID"
$ws.Range("D3").Value = "10:13 PM"
$ws.Range("E3").Value = "Leaving Camp"

# Curfew summary box values.
$ws.Range("I2").Value = "1:00 AM"
$ws.Range("I3").Value = "1:00 AM"
$ws.Range("I4").Value = "5:00 PM"

# Left/returned/still-out camp counters.
$ws.Range("I6").Value = 1
$ws.Range("I8").Value = 1

# Highlight the ID and Time In columns for the logged rows.
$ws.Range("C2").Interior.Color = 14788352
$ws.Range("C3").Interior.Color = 14788352
$ws.Range("E2").Interior.Color = 10066410
$ws.Range("E3").Interior.Color = 10066410

# Column I now holds short time values, so it no longer needs to share
# the wider width used by columns G:H.
$ws.Columns("I").ColumnWidth = 7.8
